# Policy.xlsx update:
#  - The "ID" action rule-table cell (C8) gets a new leading statement that
#    logs a debug message before building the Policy object.
#  - The row hosting that multi-line snippet (row 8) is made taller to fit
#    the extra line.
#  - The last user selection is left on D8 (the cell immediately to the
#    right of the edited one).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("policy")

# New content for C8: prepend a System.out.println(...) debug line in front
# of the existing "Policy $policy = new Policy(); $policy.setId(...)" code.
$newId = "System.out.println(""=== Fire policy!!! ==="");`r`n" + `
         "Policy `$policy = new Policy();`r`n" + `
         "`$policy.setId(""`$param"");"

$ws.Range("C8").Value = $newId

# The cell now holds three lines of text instead of two, so grow the row.
$ws.Rows.Item(8).RowHeight = 80

# Leave the selection on D8, as it was when the file was last saved.
$ws.Range("D8").Select()
